# "Darkest before the dawn"
# Sheet2 gains some source data (A1:A3) that Sheet1's new E10 formula
# references, plus a selection left on A4 (the cell right below the data).
# Sheet1 keeps the active/selected tab, so Sheet2's selection is set last
# and we re-activate Sheet1 afterwards to leave the workbook's active
# sheet unchanged.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New source values on Sheet2
$ws2.Range("A1").Value = 2
$ws2.Range("A2").Value = 3
$ws2.Range("A3").Value = 4

# New formula cell on Sheet1 referencing Sheet2
$ws1.Range("E10").Formula = "=C10+Sheet2!A2"

# Leave the selection on Sheet2 sitting at A4, then restore Sheet1 as the
# active sheet/tab (matches the original file's tabSelected on Sheet1).
$ws2.Range("A4").Select()
$ws1.Activate()
